$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range("K26").Value = 2.2
$ws.Range("Q26").Value = 2.08
$ws.Range("R26").Value = 1.73
$ws.Range("W26").Value = 2.05
$ws.Range("X26").Value = 1.7
$ws.Range("Z26").Value = 7
$ws.Range("AD26").Value = 29
$ws.Range("AE26").Value = 9
$ws.Range("AF26").Value = 7
$ws.Range("AG26").Value = 19
$ws.Range("AI26").Value = 451
$ws.Range("AP26").Value = 3.1
$ws.Range("AQ26").Value = 1.38
$ws.Range("AR26").Value = 1.6
$ws.Range("AS26").Value = 2.35

# Row 27
$ws.Range("G27").Value = 1.9
$ws.Range("H27").Value = 3.1
$ws.Range("I27").Value = 5
$ws.Range("K27").Value = 1.91
$ws.Range("U27").Value = 1.62
$ws.Range("V27").Value = 2.2
$ws.Range("X27").Value = 1.5
$ws.Range("Z27").Value = 7
$ws.Range("AF27").Value = 6.5
$ws.Range("AP27").Value = 5
$ws.Range("AR27").Value = 2.17
$ws.Range("AS27").Value = 1.7

# Row 28
$ws.Range("G28").Value = 1.6
$ws.Range("H28").Value = 3.75
$ws.Range("I28").Value = 6
$ws.Range("J28").Value = 2.25
$ws.Range("N28").Value = 9
$ws.Range("O28").Value = 1.36
$ws.Range("P28").Value = 3
$ws.Range("Q28").Value = 2.15
$ws.Range("R28").Value = 1.67
$ws.Range("S28").Value = 4
$ws.Range("T28").Value = 1.22
$ws.Range("AF28").Value = 7
$ws.Range("AK28").Value = 29
$ws.Range("AL28").Value = 19
$ws.Range("AM28").Value = 67
$ws.Range("AP28").Value = 3.1

# Row 54
$ws.Range("G54").Value = 3.6
$ws.Range("I54").Value = 2.05
$ws.Range("J54").Value = 4
$ws.Range("M54").Value = 1.05
$ws.Range("O54").Value = 1.29
$ws.Range("P54").Value = 3.5
$ws.Range("Q54").Value = 1.93
$ws.Range("R54").Value = 1.93
$ws.Range("T54").Value = 1.33
$ws.Range("Z54").Value = 19
$ws.Range("AA54").Value = 13
$ws.Range("AB54").Value = 41
$ws.Range("AC54").Value = 29
$ws.Range("AK54").Value = 10

# Row 55
$ws.Range("M55").Value = 1.11
$ws.Range("O55").Value = 1.5
$ws.Range("T55").Value = 1.17
$ws.Range("U55").Value = 1.54
$ws.Range("AQ55").Value = 1.21

# Row 71
$ws.Range("G71").Value = 2.9
$ws.Range("I71").Value = 2.55
$ws.Range("J71").Value = 3.6
$ws.Range("K71").Value = 2
$ws.Range("M71").Value = 1.08
$ws.Range("N71").Value = 8
$ws.Range("Y71").Value = 8
$ws.Range("AC71").Value = 26
$ws.Range("AJ71").Value = 7.5
$ws.Range("AM71").Value = 23

# Row 107
$ws.Range("Q107").Value = 1.93
$ws.Range("R107").Value = 1.93

# Row 108
$ws.Range("G108").Value = 2.45
$ws.Range("H108").Value = 3.6
$ws.Range("I108").Value = 2.63
$ws.Range("J108").Value = 3.1
$ws.Range("K108").Value = 2.25
$ws.Range("L108").Value = 3.25
$ws.Range("Q108").Value = 1.75
$ws.Range("R108").Value = 2.05
$ws.Range("S108").Value = 2.75
$ws.Range("T108").Value = 1.4
$ws.Range("W108").Value = 1.62
$ws.Range("X108").Value = 2.2
$ws.Range("Y108").Value = 10
$ws.Range("Z108").Value = 13
$ws.Range("AE108").Value = 13
$ws.Range("AF108").Value = 7
$ws.Range("AL108").Value = 10
$ws.Range("AM108").Value = 26
$ws.Range("AO108").Value = 26

# Row 109
$ws.Range("G109").Value = 2.15
$ws.Range("I109").Value = 3.5
$ws.Range("J109").Value = 2.75
$ws.Range("N109").Value = 10
$ws.Range("O109").Value = 1.29
$ws.Range("P109").Value = 3.5
$ws.Range("Q109").Value = 2
$ws.Range("R109").Value = 1.85
$ws.Range("W109").Value = 1.75
$ws.Range("X109").Value = 2
$ws.Range("Y109").Value = 8
$ws.Range("AD109").Value = 26
$ws.Range("AE109").Value = 10
$ws.Range("AH109").Value = 41
$ws.Range("AI109").Value = 201
$ws.Range("AJ109").Value = 11

# Row 117
$ws.Range("G117").Value = 1.65
$ws.Range("H117").Value = 3.6
$ws.Range("J117").Value = 2.38
$ws.Range("K117").Value = 2.1
$ws.Range("L117").Value = 5.5
$ws.Range("M117").Value = 1.07
$ws.Range("N117").Value = 8.5
$ws.Range("O117").Value = 1.36
$ws.Range("P117").Value = 3
$ws.Range("Q117").Value = 2.15
$ws.Range("R117").Value = 1.67
$ws.Range("S117").Value = 4
$ws.Range("T117").Value = 1.22
$ws.Range("U117").Value = 1.44
$ws.Range("V117").Value = 2.63
$ws.Range("W117").Value = 2.1
$ws.Range("X117").Value = 1.67
$ws.Range("Y117").Value = 6
$ws.Range("AB117").Value = 12
$ws.Range("AC117").Value = 15
$ws.Range("AE117").Value = 8.5
$ws.Range("AH117").Value = 67

# Row 165
$ws.Range("G165").Value = 2.35
$ws.Range("I165").Value = 2.8

# Row 166
$ws.Range("G166").Value = 1.83
$ws.Range("H166").Value = 3.1
$ws.Range("I166").Value = 4.33
$ws.Range("J166").Value = 2.63
$ws.Range("K166").Value = 1.91
$ws.Range("L166").Value = 5.5
$ws.Range("M166").Value = 1.11
$ws.Range("N166").Value = 6.5
$ws.Range("O166").Value = 1.53
$ws.Range("P166").Value = 2.38
$ws.Range("Q166").Value = 2.7
$ws.Range("R166").Value = 1.44
$ws.Range("S166").Value = 5.5
$ws.Range("T166").Value = 1.14
$ws.Range("U166").Value = 1.62
$ws.Range("V166").Value = 2.2
$ws.Range("Z166").Value = 7
$ws.Range("AB166").Value = 15
$ws.Range("AC166").Value = 21
$ws.Range("AE166").Value = 6
$ws.Range("AJ166").Value = 8.5
$ws.Range("AK166").Value = 21
$ws.Range("AL166").Value = 17

# Row 168
$ws.Range("G168").Value = 2.5
$ws.Range("H168").Value = 2.77
$ws.Range("J168").Value = 3.1
$ws.Range("K168").Value = 1.9
$ws.Range("L168").Value = 3.65
$ws.Range("O168").Value = 1.42
$ws.Range("P168").Value = 2.47
$ws.Range("Q168").Value = 2.2
$ws.Range("R168").Value = 1.52
$ws.Range("U168").Value = 1.47
$ws.Range("V168").Value = 2.32
$ws.Range("W168").Value = 1.83
$ws.Range("X168").Value = 1.78
$ws.Range("Y168").Value = 7
$ws.Range("Z168").Value = 12
$ws.Range("AB168").Value = 29
$ws.Range("AC168").Value = 23
$ws.Range("AD168").Value = 35
$ws.Range("AE168").Value = 6.8
$ws.Range("AF168").Value = 5.5
$ws.Range("AG168").Value = 14.5
$ws.Range("AH168").Value = 80
$ws.Range("AI168").Value = 700
$ws.Range("AJ168").Value = 7.5
$ws.Range("AM168").Value = 40
